$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.498.17"
$c.Style = $s

$c = $ws.Range("E2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.05%  "
$c.Style = $s

$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.864.61"
$c.Style = $s

$c = $ws.Range("E3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.58%  "
$c.Style = $s

$c = $ws.Range("E4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.48%  "
$c.Style = $s

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "315.81"
$c.Style = $s

$c = $ws.Range("E5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.22%  "
$c.Style = $s

$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = $s

$c = $ws.Range("E6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.44%  "
$c.Style = $s

$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4678"
$c.Style = $s

$c = $ws.Range("E7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.79%  "
$c.Style = $s

$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3732"
$c.Style = $s

$c = $ws.Range("E8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.23%  "
$c.Style = $s

$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07382"
$c.Style = $s

$c = $ws.Range("E9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.24%  "
$c.Style = $s

$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8893"
$c.Style = $s

$c = $ws.Range("E10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.54%  "
$c.Style = $s

$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07952"
$c.Style = $s

$c = $ws.Range("E11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.12%  "
$c.Style = $s

$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.03"
$c.Style = $s

$c = $ws.Range("E12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.43%  "
$c.Style = $s

$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.859.83"
$c.Style = $s

$c = $ws.Range("E13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +8.13%  "
$c.Style = $s

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.426"
$c.Style = $s

$c = $ws.Range("E14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.03%  "
$c.Style = $s

$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.611"
$c.Style = $s

$c = $ws.Range("E15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.91%  "
$c.Style = $s

$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "92.78"
$c.Style = $s

$c = $ws.Range("E16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.10%  "
$c.Style = $s

$c = $ws.Range("E17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.33%  "
$c.Style = $s

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000008946"
$c.Style = $s

$c = $ws.Range("E19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.43%  "
$c.Style = $s

$c = $ws.Range("E20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.42%  "
$c.Style = $s

$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.526.85"
$c.Style = $s

$c = $ws.Range("E21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.34%  "
$c.Style = $s

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.158"
$c.Style = $s

$c = $ws.Range("E22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.Style = $s

$c = $ws.Range("E23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.64%  "
$c.Style = $s

$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.100.46"
$c.Style = $s

$c = $ws.Range("E24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +4.01%  "
$c.Style = $s

$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "153.32"
$c.Style = $s

$c = $ws.Range("E25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.96%  "
$c.Style = $s

$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.873"
$c.Style = $s

$c = $ws.Range("E26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.22%  "
$c.Style = $s

$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.53"
$c.Style = $s

$c = $ws.Range("E27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.51%  "
$c.Style = $s

$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.091"
$c.Style = $s

$c = $ws.Range("E28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.19%  "
$c.Style = $s

$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.169"
$c.Style = $s

$c = $ws.Range("E29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.27%  "
$c.Style = $s

$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "117.16"
$c.Style = $s

$c = $ws.Range("E30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.67%  "
$c.Style = $s

$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.08917"
$c.Style = $s

$c = $ws.Range("E31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.59%  "
$c.Style = $s

$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7556"
$c.Style = $s

$c = $ws.Range("E32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.89%  "
$c.Style = $s

$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.021"
$c.Style = $s

$c = $ws.Range("E33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.Style = $s

$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.166"
$c.Style = $s

$c = $ws.Range("E34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.48%  "
$c.Style = $s

$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.491"
$c.Style = $s

$c = $ws.Range("E35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.Style = $s

$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.643"
$c.Style = $s

$c = $ws.Range("E36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +9.59%  "
$c.Style = $s

$c = $ws.Range("B37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "TrustWalletToken"
$c.Style = $s

$c = $ws.Range("C37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.Style = $s

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.084"
$c.Style = $s

$c = $ws.Range("E37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.83%  "
$c.Style = $s

$c = $ws.Range("B38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.Style = $s

$c = $ws.Range("C38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = $s

$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01971"
$c.Style = $s

$c = $ws.Range("E38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.61%  "
$c.Style = $s

$c = $ws.Range("E39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.85%  "
$c.Style = $s

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.991"
$c.Style = $s

$c = $ws.Range("E40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.57%  "
$c.Style = $s

$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.194"
$c.Style = $s

$c = $ws.Range("E41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.85%  "
$c.Style = $s

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5221"
$c.Style = $s

$c = $ws.Range("E42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.53%  "
$c.Style = $s

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1647"
$c.Style = $s

$c = $ws.Range("E43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.54%  "
$c.Style = $s

$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.367"
$c.Style = $s

$c = $ws.Range("E44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.68%  "
$c.Style = $s

$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4882"
$c.Style = $s

$c = $ws.Range("E45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.66%  "
$c.Style = $s

$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.39"
$c.Style = $s

$c = $ws.Range("E46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.45%  "
$c.Style = $s

$c = $ws.Range("E47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.48%  "
$c.Style = $s

$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "103.87"
$c.Style = $s

$c = $ws.Range("E48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.99%  "
$c.Style = $s

$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.662"
$c.Style = $s

$c = $ws.Range("E49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.09%  "
$c.Style = $s

$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06261"
$c.Style = $s

$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "66.00"
$c.Style = $s

$c = $ws.Range("E51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.94%  "
$c.Style = $s

